$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Food"
$ws.Range("D1").Value = "Wood"
$ws.Range("E1").Value = "Stone"
$ws.Range("F1").Value = "Metal"
$ws.Range("G1").Value = "Pollution"
$ws.Range("H1").Value = "Cost"

# Header row used to carry a "Hyperlink-ish" style on C1 and E1 (no visible
# hyperlink, just the cell style / font left over). In the new layout that
# style moved to G1 (Pollution) instead, with a plain (no underline, no
# special color) font.
$ws.Range("C1").Style = "Normal"
$ws.Range("E1").Style = "Normal"
$ws.Range("G1").Style = "Hyperlink"
$ws.Range("G1").Font.Underline = $false
$ws.Range("G1").Font.ThemeColor = 1

# --- Row 2 : Town Center ---
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = 200
$ws.Range("E2").Value = 200
$ws.Range("F2").Value = 300
$ws.Range("G2").Value = "pol_prog"
$ws.Range("H2").Value = "cost_prog"

# --- Row 3 : Farm ---
$ws.Range("C3").Value = 200
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = "pol_prog"
$ws.Range("H3").Value = "cost_prog"

# --- Row 4 : Factory ---
$ws.Range("C4").Value = 300
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 150
$ws.Range("G4").Value = "pol_prog"
$ws.Range("H4").Value = "cost_prog"

# --- Row 5 : Filteration Plant ---
$ws.Range("C5").Value = 400
$ws.Range("D5").Value = 300
$ws.Range("E5").Value = 200
$ws.Range("F5").Value = 100
$ws.Range("G5").Value = "pol_prog"
$ws.Range("H5").Value = "cost_prog"

# --- Row 6 : House ---
$ws.Range("C6").Value = 150
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = "pol_prog"
$ws.Range("H6").Value = "cost_prog"

# --- Column widths for the new columns ---
# (target stored widths are 9.42578125 / 12.140625 / 12 ; the runtime quantizes
# ColumnWidth to 1/6-character increments, so these are the closest reachable
# values under that rounding)
$ws.Columns.Item(6).ColumnWidth = 8.666666666666666
$ws.Columns.Item(7).ColumnWidth = 11.333333333333334
$ws.Columns.Item(8).ColumnWidth = 11.166666666666666

# --- Selection moves to F2 ---
$ws.Range("F2").Select() | Out-Null

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1
